# check_if_jails_days_imposed greater than suspended — append case 21TRD09437
# (Hemmeter) charges, and drop the stray empty G1093 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Row 1093 had an empty placeholder cell in column G; clear it so the
#        cell disappears entirely (matches the diff removing <c r="G1093".../>).
$ws.Range("G1093").Value = ""

# --- 2. Append 12 new rows (1095-1106) repeating the case's charge pattern.
#        Columns H/I/J/K (and the numeric-looking D/J/K entries) must stay
#        TEXT, not auto-converted to numbers/currency, so pre-format those
#        cells as Text before assigning their values.

$newRows = @(
    @{ A="21TRD09437"; B="Hemmeter"; C="DUS "; D="4510.11"; E="M1"; F="No Contest"; G="Guilty"; H="$ 0"; I="$ 0"; J="50"; K="5" },
    @{ A="21TRD09437"; B="Hemmeter"; C="1st Speed 1 Yr School >35mphm4"; D="4511.21B1A"; E="M4"; F="No Contest"; G="Guilty"; H="$ 0"; I="$ 0"; J="10"; K="20" },
    @{ A="21TRD09437"; B="Hemmeter"; C="Reckless Operation 1st In 1 Yr"; D="4511.20"; E="MM"; F="Dismissed"; G=$null; H=" "; I=" "; J=" "; K=" " },
    @{ A="21TRD09437"; B="Hemmeter"; C="DUS "; D="4510.11"; E="M1"; F="No Contest"; G="Guilty"; H="$ 0"; I="$ 0"; J="50"; K="5" },
    @{ A="21TRD09437"; B="Hemmeter"; C="1st Speed 1 Yr School >35mphm4"; D="4511.21B1A"; E="M4"; F="No Contest"; G="Guilty"; H="$ 0"; I="$ 0"; J="10"; K="20" },
    @{ A="21TRD09437"; B="Hemmeter"; C="Reckless Operation 1st In 1 Yr"; D="4511.20"; E="MM"; F="Dismissed"; G=$null; H=" "; I=" "; J=" "; K=" " },
    @{ A="21TRD09437"; B="Hemmeter"; C="DUS "; D="4510.11"; E="M1"; F="No Contest"; G="Guilty"; H="$ 0"; I="$ 0"; J="50"; K="5" },
    @{ A="21TRD09437"; B="Hemmeter"; C="1st Speed 1 Yr School >35mphm4"; D="4511.21B1A"; E="M4"; F="No Contest"; G="Guilty"; H="$ 0"; I="$ 0"; J="10"; K="20" },
    @{ A="21TRD09437"; B="Hemmeter"; C="Reckless Operation 1st In 1 Yr"; D="4511.20"; E="MM"; F="Dismissed"; G=$null; H=" "; I=" "; J=" "; K=" " },
    @{ A="21TRD09437"; B="Hemmeter"; C="DUS "; D="4510.11"; E="M1"; F="No Contest"; G="Guilty"; H="$ 0"; I="$ 0"; J="50"; K="5" },
    @{ A="21TRD09437"; B="Hemmeter"; C="1st Speed 1 Yr School >35mphm4"; D="4511.21B1A"; E="M4"; F="No Contest"; G="Guilty"; H="$ 0"; I="$ 0"; J="10"; K="20" },
    @{ A="21TRD09437"; B="Hemmeter"; C="Reckless Operation 1st In 1 Yr"; D="4511.20"; E="MM"; F="Dismissed"; G=""; H=" "; I=" "; J=" "; K=" " }
)

# Column D/J/K values that are pure numbers/decimals, plus every H/I "$ 0",
# must be pre-formatted as Text so the COM layer doesn't coerce them into
# numeric/currency cells.
$textForceCols = @("D", "H", "I", "J", "K")

$startRow = 1095
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K")) {
        $val = $row[$col]
        if ($null -eq $val) {
            continue
        }
        $cell = $ws.Range("$col$r")
        if ($textForceCols -contains $col) {
            $cell.NumberFormat = "@"
            $cell.Value = $val
            # Drop the Text number-format again now that the value is safely
            # stored as text, so we don't leave a stray style behind.
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }
}
